# Adds the "2022-Q3" quarterly holdings sheet:
#  1. Inserts a new summary row (2022-Q3) at the top of the "总计" sheet,
#     pushing the existing quarters down by one row.
#  2. Inserts a brand-new worksheet named "2022-Q3" right after "总计",
#     containing the per-fund holding detail for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet - insert new row 2 for 2022-Q3
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push existing data rows (old row 2..8) down by one, duplicating the
# formatting of the row being displaced (Excel's native Insert behaviour).
$summary.Rows.Item(2).Insert()

# The freshly inserted row borrowed formatting from row 1 (the header),
# which has no entries under columns B:D -> strip that back to the default
# "Normal" style so the new data row matches its sibling data rows exactly.
$summary.Range("B2:D2").Style = "Normal"

# Fill in the new 2022-Q3 summary row.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 8.01

# Column A holds a simple 0-based row counter; recompute it for every data
# row now that one more row exists (0,1,2,...,7).
$labels = @("2022-Q3", "2022-Q2", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = 2 + $i
    $summary.Range("A$r").Value = $i
}
# Column A keeps the bold/centered/bordered style used throughout the sheet;
# copy it from the (untouched) row below onto the newly numbered A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("A2").Value = 0

# ---------------------------------------------------------------------------
# 2. Brand-new "2022-Q3" worksheet with the per-fund holdings detail
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$q3.Name = "2022-Q3"

# Borrow the header-row formatting from the sibling "2022-Q2" sheet
# (bold/border header cells in B1:H1) so the new sheet matches its
# siblings exactly.
# (Re-fetch the "2022-Q2" sheet by name instead of reusing the handle that
# was passed into Add() above - that handle no longer targets live cells.)
$afterQ2 = $wb.Worksheets.Item("2022-Q2")
$afterQ2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Every row's column-A index cell carries the same bold/border/centered
# style; copy it from a single template cell onto the whole A2:A13 block
# in one shot (Excel tiles a 1x1 source across a larger paste target).
$afterQ2.Range("A2").Copy()
$q3.Range("A2:A13").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $col = [char]([int][char]'B' + $c)
    $addr = "$col" + "1"
    # B1:H1 already carry the bold/border header style copied above; none of
    # these Chinese labels look numeric, so a plain Value assignment can't
    # be mis-coerced and the pre-existing "s=2" style is left untouched.
    $q3.Range($addr).Value = $headers[$c]
}

$rows = @(
    @("001257", "兴业收益增强债券A",               "68.35", "20.08", "3.36", "2.2966", 2),
    @("006567", "中泰星元价值优选灵活配置混合A",    "46.48", "87.07", "3.50", "1.6268", 10),
    @("005984", "兴业聚华混合A",                    "23.61", "29.42", "4.83", "1.1404", 1),
    @("006624", "中泰玉衡价值优选混合A",             "23.07", "88.61", "3.50", "0.8074", 9),
    @("012940", "中泰星元价值优选灵活配置混合C",    "22.67", "87.07", "3.50", "0.7934", 10),
    @("000893", "工银创新动力股票",                 "11.92", "81.96", "4.32", "0.5149", 3),
    @("001258", "兴业收益增强债券C",                "11.01", "20.08", "3.36", "0.3699", 2),
    @("005985", "兴业聚华混合C",                     "6.45", "29.42", "4.83", "0.3115", 1),
    @("002076", "浙商中证500指数增强A",              "6.79", "85.72", "1.08", "0.0733", 7),
    @("016090", "中泰玉衡价值优选混合C",              "1.01", "88.61", "3.50", "0.0354", 9),
    @("007386", "浙商中证500指数增强C",               "1.97", "85.72", "1.08", "0.0213", 7),
    @("510660", "华夏上证医药卫生ETF",                "0.96", "99.53", "1.99", "0.0191", 10)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    # Column A keeps the bold/border/centered "s=2" style copied above
    # (from the template's A2); plain numbers can't be mis-coerced so this
    # plain assignment leaves that style untouched.
    $q3.Range("A$r").Value = $i

    # Text-like columns (fund code / name / size / position / ratio / value)
    # must stay text even though several look numeric (e.g. leading-zero
    # fund codes, "3.50"), matching the source data's inline-string typing.
    # These columns carry no explicit style in the source, so after forcing
    # Text entry we strip the "Text" number-format style back off again.
    $cols = @("B", "C", "D", "E", "F", "G")
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = "$($cols[$c])$r"
        $q3.Range($addr).NumberFormat = "@"
        $q3.Range($addr).Value = $row[$c]
        $q3.Range($addr).ClearFormats()
    }

    # Rank column is a genuine number, unstyled.
    $q3.Range("H$r").Value = $row[6]
}
